$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (engine stores width = ColumnWidth + 5/6, so back out 36)
$ws.Columns.Item(1).ColumnWidth = 35.166666666666664

# Add new row of data (set B11 first so shared-string order matches)
$ws.Range("B11").Value = "https://youtu.be/oC9N6lz70kY"
$ws.Range("A11").Value = "Практические задания 3, 4, 5 и 6"

# Add hyperlink on B11 and apply hyperlink style
$ws.Hyperlinks.Add($ws.Range("B11"), "https://youtu.be/oC9N6lz70kY")
$ws.Range("B11").Style = "Гиперссылка"

# Update selection to reflect last active cell
$ws.Range("K11").Select()
